# Scheduled data refresh: update market-price / profit figures on each class sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1067.3077
$ws.Range("I2").Value = 541.05884
$ws.Range("J2").Value = 2061.3333
$ws.Range("K2").Value = 541.05884
$ws.Range("L2").Value = 2061.3333
$ws.Range("M2").Value = -428.05884
$ws.Range("N2").Value = -2287.3333
$ws.Range("H4").Value = 166.66667
$ws.Range("I4").Value = 166.66667
$ws.Range("K4").Value = 166.66667
$ws.Range("M4").Value = -52.66667000000001
$ws.Range("I32").Value = 333333760
$ws.Range("J32").Value = 797.6667
$ws.Range("K32").Value = 333333760
$ws.Range("L32").Value = 797.6667
$ws.Range("M32").Value = -333333434
$ws.Range("N32").Value = -1449.6667
$ws.Range("H51").Value = 1800.125
$ws.Range("I51").Value = 1170
$ws.Range("J51").Value = 2010.1666
$ws.Range("K51").Value = 1170
$ws.Range("L51").Value = 2010.1666
$ws.Range("M51").Value = -686
$ws.Range("N51").Value = -2978.1666
$ws.Range("H53").Value = 143.63158
$ws.Range("J53").Value = 276.5
$ws.Range("L53").Value = 276.5
$ws.Range("N53").Value = -1550.5
$ws.Range("H62").Value = 11718.8
$ws.Range("I62").Value = 1909.7778
$ws.Range("K62").Value = 1909.7778
$ws.Range("M62").Value = -1285.7778
$ws.Range("H65").Value = 11718.8
$ws.Range("I65").Value = 1909.7778
$ws.Range("K65").Value = 9548.889000000001
$ws.Range("M65").Value = -6428.889000000001
$ws.Range("H98").Value = 2538
$ws.Range("I98").Value = 1190.75
$ws.Range("K98").Value = 1190.75
$ws.Range("M98").Value = 307.25
$ws.Range("H116").Value = 4588
$ws.Range("I116").Value = 2058.1333
$ws.Range("J116").Value = 7507.077
$ws.Range("K116").Value = 2058.1333
$ws.Range("L116").Value = 7507.077
$ws.Range("M116").Value = 1383.8667
$ws.Range("N116").Value = -14391.077
$ws.Range("H122").Value = 2538
$ws.Range("I122").Value = 1190.75
$ws.Range("K122").Value = 3572.25
$ws.Range("M122").Value = -1122.25
$ws.Range("H137").Value = 2663.3157
$ws.Range("I137").Value = 2183.3333
$ws.Range("J137").Value = 2884.8462
$ws.Range("K137").Value = 6549.999899999999
$ws.Range("L137").Value = 8654.5386
$ws.Range("M137").Value = -3999.999899999999
$ws.Range("N137").Value = -13754.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2039.2
$ws.Range("I88").Value = 1932
$ws.Range("K88").Value = 1932
$ws.Range("M88").Value = -1526
$ws.Range("H91").Value = 2039.2
$ws.Range("I91").Value = 1932
$ws.Range("K91").Value = 1932
$ws.Range("M91").Value = -528
$ws.Range("H132").Value = 1780.909
$ws.Range("I132").Value = 1501.1818
$ws.Range("J132").Value = 2620.0908
$ws.Range("K132").Value = 4503.5454
$ws.Range("L132").Value = 7860.2724
$ws.Range("M132").Value = -1973.5454
$ws.Range("N132").Value = -12920.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 60000
$ws.Range("J137").Value = 60000
$ws.Range("L137").Value = 60000
$ws.Range("N137").Value = -70200
$ws.Range("H140").Value = 55650.715
$ws.Range("J140").Value = 55650.715
$ws.Range("L140").Value = 55650.715
$ws.Range("N140").Value = -66010.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2255.0222
$ws.Range("I31").Value = 1449.55
$ws.Range("J31").Value = 8698.799999999999
$ws.Range("K31").Value = 1449.55
$ws.Range("L31").Value = 8698.799999999999
$ws.Range("M31").Value = -1154.55
$ws.Range("N31").Value = -9288.799999999999
$ws.Range("H34").Value = 2255.0222
$ws.Range("I34").Value = 1449.55
$ws.Range("J34").Value = 8698.799999999999
$ws.Range("K34").Value = 1449.55
$ws.Range("L34").Value = 8698.799999999999
$ws.Range("M34").Value = -1247.55
$ws.Range("N34").Value = -9102.799999999999
$ws.Range("H69").Value = 8475.125
$ws.Range("I69").Value = 4267
$ws.Range("K69").Value = 4267
$ws.Range("M69").Value = -3518
$ws.Range("H72").Value = 8475.125
$ws.Range("I72").Value = 4267
$ws.Range("K72").Value = 12801
$ws.Range("M72").Value = -9057
$ws.Range("H87").Value = 23500
$ws.Range("J87").Value = 23500
$ws.Range("L87").Value = 23500
$ws.Range("N87").Value = -25872
$ws.Range("H90").Value = 23500
$ws.Range("J90").Value = 23500
$ws.Range("L90").Value = 70500
$ws.Range("N90").Value = -82356

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 350.46667
$ws.Range("I40").Value = 61.416668
$ws.Range("J40").Value = 1506.6666
$ws.Range("K40").Value = 245.666672
$ws.Range("L40").Value = 6026.6664
$ws.Range("M40").Value = -176.666672
$ws.Range("N40").Value = -6164.6664
$ws.Range("H92").Value = 199
$ws.Range("I92").Value = 199
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 597
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("M92").Value = 651

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 15700.728
$ws.Range("J123").Value = 15700.728
$ws.Range("L123").Value = 15700.728
$ws.Range("N123").Value = -20600.728

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5885799.5
$ws.Range("I7").Value = 11113467
$ws.Range("J7").Value = 4673.875
$ws.Range("K7").Value = 11113467
$ws.Range("L7").Value = 4673.875
$ws.Range("M7").Value = -11113355
$ws.Range("N7").Value = -4897.875
$ws.Range("H22").Value = 859.7727
$ws.Range("I22").Value = 769
$ws.Range("K22").Value = 769
$ws.Range("M22").Value = -474
$ws.Range("H27").Value = 859.7727
$ws.Range("I27").Value = 769
$ws.Range("K27").Value = 769
$ws.Range("M27").Value = -662
$ws.Range("H126").Value = 5885799.5
$ws.Range("I126").Value = 11113467
$ws.Range("J126").Value = 4673.875
$ws.Range("K126").Value = 33340401
$ws.Range("L126").Value = 14021.625
$ws.Range("M126").Value = -33337931
$ws.Range("N126").Value = -18961.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6383.9
$ws.Range("J81").Value = 1617.3636
$ws.Range("L81").Value = 3234.7272
$ws.Range("N81").Value = -5356.727199999999
$ws.Range("H84").Value = 6383.9
$ws.Range("J84").Value = 1617.3636
$ws.Range("L84").Value = 16173.636
$ws.Range("N84").Value = -26781.636
$ws.Range("H123").Value = 43944
$ws.Range("J123").Value = 43944
$ws.Range("L123").Value = 43944
$ws.Range("N123").Value = -53744
